$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: 2021年 ---
$ws.Range("A7").Value = "2021年"
$ws.Range("B7").Value = 102.3
$ws.Range("C7").Value = 100.8
$ws.Range("D7").Value = 101.8
$ws.Range("E7").Value = 100.4

# --- Row 8: 2022年 ---
$ws.Range("A8").Value = "2022年"
$ws.Range("C8").Value = 100.7

# The source sheet leaves some trailing cells present but blank (inline empty
# strings), e.g. F5/F6, B8/D8/E8/F8, F7. Replicate that by pasting the
# existing blank cell F6 (already an empty-string cell) into those spots so
# they materialize as real (empty) cells instead of being entirely absent.
$ws.Range("F6").Copy() | Out-Null
$ws.Range("F7:F8").PasteSpecial(-4122) | Out-Null
$ws.Range("B8").PasteSpecial(-4122) | Out-Null
$ws.Range("D8").PasteSpecial(-4122) | Out-Null
$ws.Range("E8").PasteSpecial(-4122) | Out-Null

# Re-apply the year-label formatting (bold, centered, bordered) used by the
# rest of column A to the two newly added year cells.
$ws.Range("A6").Copy() | Out-Null
$ws.Range("A7:A8").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0
